# Update the selection on the "category" sheet (F7 -> F9) without leaving it as the active tab.
$wb = $excel.ActiveWorkbook
$catWs = $wb.Worksheets.Item("category")
$catWs.Activate()
$catWs.Range("F9").Select()

# Rename "listContent" -> "contentCategory" and make it the active sheet.
$contentCategoryWs = $wb.Worksheets.Item("listContent")
$contentCategoryWs.Name = "contentCategory"
$contentCategoryWs.Activate()
